$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 49 (existing rows 49-87 shift down to become 52-90)
$ws.Rows("49:51").Insert()

# New weekly entries (Feria Lagunitas de Puerto Montt, Cereza) for 2022-12-09
# (Excel serial date 44904 == 2022-12-09)
$newDate = 44904

$rows = @(
  @{ Row=49; Variedad="Brooks";  Calidad="Primera"; Volumen=800; PMin=8000;  PMax=8500;  PProm=8250;  Unidad="$/bandeja 10 kilos"; Origen="Provincia de Curicó"; PKg=825;  KgUnidad=10 },
  @{ Row=50; Variedad="Lapins";  Calidad="Primera"; Volumen=800; PMin=8000;  PMax=8500;  PProm=8250;  Unidad="$/bandeja 10 kilos"; Origen="Provincia de Curicó"; PKg=825;  KgUnidad=10 },
  @{ Row=51; Variedad="Rainier"; Calidad="Primera"; Volumen=600; PMin=10000; PMax=11000; PProm=10500; Unidad="$/bandeja 10 kilos"; Origen="Provincia de Curicó"; PKg=1050; KgUnidad=10 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 4
    $ws.Cells.Item($row, 2).Value = "Feria Lagunitas de Puerto Montt"
    $ws.Cells.Item($row, 3).Value = "Los Lagos"
    $ws.Cells.Item($row, 4).Value = $newDate
    $ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 5).Value = 10
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100103
    $ws.Cells.Item($row, 8).Value = "Frutos de hueso (carozo)"
    $ws.Cells.Item($row, 9).Value = 100103001
    $ws.Cells.Item($row, 10).Value = "Cereza"
    $ws.Cells.Item($row, 11).Value = $r.Variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.PMin
    $ws.Cells.Item($row, 15).Value = $r.PMax
    $ws.Cells.Item($row, 16).Value = $r.PProm
    $ws.Cells.Item($row, 17).Value = $r.Unidad
    $ws.Cells.Item($row, 18).Value = $r.Origen
    $ws.Cells.Item($row, 19).Value = $r.PKg
    $ws.Cells.Item($row, 20).Value = $r.KgUnidad
}
